# Insert a new row at position 592 (this shifts the existing rows 592:633
# down to 593:634 and extends the sheet's used range to row 634), then
# populate the new row:
#   A592 = 2026/01/10   B592 = 土   C592 = 10   D592 = 25
# A592/B592 are filled by copying the date/weekday from the row above
# (which already carries the same "2026/01/10" / "土" labels) so the new
# cells come out as plain text, matching the rest of the column, instead
# of being auto-converted to a date serial number by a literal-text
# assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A592:D592").EntireRow.Insert()

$ws.Range("A591:B591").Copy()
$ws.Range("A592").PasteSpecial(-4163)

$ws.Range("C592").Value = 10
$ws.Range("D592").Value = 25
